$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 47-92: "freq" dataset entries appended to the evolution log

# Row 47
$ws.Range("A47").Value = 1753871814
$ws.Range("B47").Value = "add"
$ws.Range("C47").Value = "freq"
$ws.Range("D47").Value = "variable_1---null"

# Row 48
$ws.Range("A48").Value = 1753871814
$ws.Range("B48").Value = "add"
$ws.Range("C48").Value = "freq"
$ws.Range("D48").Value = "variable_1---value_1"

# Row 49
$ws.Range("A49").Value = 1753871814
$ws.Range("B49").Value = "add"
$ws.Range("C49").Value = "freq"
$ws.Range("D49").Value = "variable_1---value_2"

# Row 50
$ws.Range("A50").Value = 1753871964
$ws.Range("B50").Value = "add"
$ws.Range("C50").Value = "freq"
$ws.Range("D50").Value = "accident_route___variable_1---value"

# Row 51
$ws.Range("A51").Value = 1753871964
$ws.Range("B51").Value = "add"
$ws.Range("C51").Value = "freq"
$ws.Range("D51").Value = "accident_route___variable_1---value_1"

# Row 52
$ws.Range("A52").Value = 1753871964
$ws.Range("B52").Value = "add"
$ws.Range("C52").Value = "freq"
$ws.Range("D52").Value = "accident_route___variable_1---value_2"

# Row 53
$ws.Range("A53").Value = 1753871964
$ws.Range("B53").Value = "add"
$ws.Range("C53").Value = "freq"
$ws.Range("D53").Value = "accident_route___variable_1---value_3"

# Row 54
$ws.Range("A54").Value = 1753871964
$ws.Range("B54").Value = "delete"
$ws.Range("C54").Value = "freq"
$ws.Range("D54").Value = "variable_1---null"
$ws.Range("E54").Value = "variable_1"

# Row 55
$ws.Range("A55").Value = 1753871964
$ws.Range("B55").Value = "delete"
$ws.Range("C55").Value = "freq"
$ws.Range("D55").Value = "variable_1---value_1"
$ws.Range("E55").Value = "variable_1"

# Row 56
$ws.Range("A56").Value = 1753871964
$ws.Range("B56").Value = "delete"
$ws.Range("C56").Value = "freq"
$ws.Range("D56").Value = "variable_1---value_2"
$ws.Range("E56").Value = "variable_1"

# Row 57
$ws.Range("A57").Value = 1753874536
$ws.Range("B57").Value = "add"
$ws.Range("C57").Value = "freq"
$ws.Range("D57").Value = "accident_route___variable_13---value"

# Row 58
$ws.Range("A58").Value = 1753874536
$ws.Range("B58").Value = "add"
$ws.Range("C58").Value = "freq"
$ws.Range("D58").Value = "accident_route___variable_13---value_1"

# Row 59
$ws.Range("A59").Value = 1753874536
$ws.Range("B59").Value = "add"
$ws.Range("C59").Value = "freq"
$ws.Range("D59").Value = "accident_route___variable_13---value_2"

# Row 60
$ws.Range("A60").Value = 1753874536
$ws.Range("B60").Value = "add"
$ws.Range("C60").Value = "freq"
$ws.Range("D60").Value = "accident_route___variable_13---value_3"

# Row 61
$ws.Range("A61").Value = 1753874536
$ws.Range("B61").Value = "delete"
$ws.Range("C61").Value = "freq"
$ws.Range("D61").Value = "accident_route___variable_1---value"
$ws.Range("E61").Value = "accident_route___variable_1"

# Row 62
$ws.Range("A62").Value = 1753874536
$ws.Range("B62").Value = "delete"
$ws.Range("C62").Value = "freq"
$ws.Range("D62").Value = "accident_route___variable_1---value_1"
$ws.Range("E62").Value = "accident_route___variable_1"

# Row 63
$ws.Range("A63").Value = 1753874536
$ws.Range("B63").Value = "delete"
$ws.Range("C63").Value = "freq"
$ws.Range("D63").Value = "accident_route___variable_1---value_2"
$ws.Range("E63").Value = "accident_route___variable_1"

# Row 64
$ws.Range("A64").Value = 1753874536
$ws.Range("B64").Value = "delete"
$ws.Range("C64").Value = "freq"
$ws.Range("D64").Value = "accident_route___variable_1---value_3"
$ws.Range("E64").Value = "accident_route___variable_1"

# Row 65
$ws.Range("A65").Value = 1753874564
$ws.Range("B65").Value = "add"
$ws.Range("C65").Value = "freq"
$ws.Range("D65").Value = "ser_pub_loc___variable_13---value"

# Row 66
$ws.Range("A66").Value = 1753874564
$ws.Range("B66").Value = "add"
$ws.Range("C66").Value = "freq"
$ws.Range("D66").Value = "ser_pub_loc___variable_13---value_1"

# Row 67
$ws.Range("A67").Value = 1753874564
$ws.Range("B67").Value = "add"
$ws.Range("C67").Value = "freq"
$ws.Range("D67").Value = "ser_pub_loc___variable_13---value_2"

# Row 68
$ws.Range("A68").Value = 1753874564
$ws.Range("B68").Value = "add"
$ws.Range("C68").Value = "freq"
$ws.Range("D68").Value = "ser_pub_loc___variable_13---value_3"

# Row 69
$ws.Range("A69").Value = 1753874564
$ws.Range("B69").Value = "delete"
$ws.Range("C69").Value = "freq"
$ws.Range("D69").Value = "accident_route___variable_13---value"
$ws.Range("E69").Value = "accident_route___variable_13"

# Row 70
$ws.Range("A70").Value = 1753874564
$ws.Range("B70").Value = "delete"
$ws.Range("C70").Value = "freq"
$ws.Range("D70").Value = "accident_route___variable_13---value_1"
$ws.Range("E70").Value = "accident_route___variable_13"

# Row 71
$ws.Range("A71").Value = 1753874564
$ws.Range("B71").Value = "delete"
$ws.Range("C71").Value = "freq"
$ws.Range("D71").Value = "accident_route___variable_13---value_2"
$ws.Range("E71").Value = "accident_route___variable_13"

# Row 72
$ws.Range("A72").Value = 1753874564
$ws.Range("B72").Value = "delete"
$ws.Range("C72").Value = "freq"
$ws.Range("D72").Value = "accident_route___variable_13---value_3"
$ws.Range("E72").Value = "accident_route___variable_13"

# Row 73
$ws.Range("A73").Value = 1753875837
$ws.Range("B73").Value = "add"
$ws.Range("C73").Value = "freq"
$ws.Range("D73").Value = "ser_pub_loc___variable_13---value_4"

# Row 74
$ws.Range("A74").Value = 1753875837
$ws.Range("B74").Value = "add"
$ws.Range("C74").Value = "freq"
$ws.Range("D74").Value = "ser_pub_loc___variable_13---value_5"

# Row 75
$ws.Range("A75").Value = 1753875837
$ws.Range("B75").Value = "add"
$ws.Range("C75").Value = "freq"
$ws.Range("D75").Value = "ser_pub_loc___variable_13---value_6"

# Row 76
$ws.Range("A76").Value = 1753875837
$ws.Range("B76").Value = "add"
$ws.Range("C76").Value = "freq"
$ws.Range("D76").Value = "ser_pub_loc___variable_13---value_7"

# Row 77
$ws.Range("A77").Value = 1753875837
$ws.Range("B77").Value = "add"
$ws.Range("C77").Value = "freq"
$ws.Range("D77").Value = "ser_pub_loc___variable_13---value_8"

# Row 78
$ws.Range("A78").Value = 1753875837
$ws.Range("B78").Value = "add"
$ws.Range("C78").Value = "freq"
$ws.Range("D78").Value = "ser_pub_loc___variable_13---value_9"

# Row 79
$ws.Range("A79").Value = 1753875837
$ws.Range("B79").Value = "add"
$ws.Range("C79").Value = "freq"
$ws.Range("D79").Value = "ser_pub_loc___variable_13---value_10"

# Row 80
$ws.Range("A80").Value = 1753875837
$ws.Range("B80").Value = "add"
$ws.Range("C80").Value = "freq"
$ws.Range("D80").Value = "ser_pub_loc___variable_13---value_11"

# Row 81
$ws.Range("A81").Value = 1753875837
$ws.Range("B81").Value = "add"
$ws.Range("C81").Value = "freq"
$ws.Range("D81").Value = "ser_pub_loc___variable_13---value_12"

# Row 82
$ws.Range("A82").Value = 1753875837
$ws.Range("B82").Value = "add"
$ws.Range("C82").Value = "freq"
$ws.Range("D82").Value = "ser_pub_loc___variable_13---value_13"

# Row 83
$ws.Range("A83").Value = 1753875837
$ws.Range("B83").Value = "add"
$ws.Range("C83").Value = "freq"
$ws.Range("D83").Value = "ser_pub_loc___variable_13---value_14"

# Row 84
$ws.Range("A84").Value = 1753875837
$ws.Range("B84").Value = "update"
$ws.Range("C84").Value = "freq"
$ws.Range("D84").Value = "ser_pub_loc___variable_13---value"
$ws.Range("F84").Value = "freq"
$ws.Range("G84").NumberFormat = "@"
$ws.Range("G84").Value = "200"
$ws.Range("H84").NumberFormat = "@"
$ws.Range("H84").Value = "50"

# Row 85
$ws.Range("A85").Value = 1753876042
$ws.Range("B85").Value = "add"
$ws.Range("C85").Value = "freq"
$ws.Range("D85").Value = "ser_pub_loc___variable_13---value_2 long with a lot of text"

# Row 86
$ws.Range("A86").Value = 1753876042
$ws.Range("B86").Value = "delete"
$ws.Range("C86").Value = "freq"
$ws.Range("D86").Value = "ser_pub_loc___variable_13---value_2"
$ws.Range("E86").Value = "ser_pub_loc___variable_13"

# Row 87
$ws.Range("A87").Value = 1753876086
$ws.Range("B87").Value = "add"
$ws.Range("C87").Value = "freq"
$ws.Range("D87").Value = "ser_pub_loc___variable_13---value_2 long with a lot of text, but like a lot of text very very very"

# Row 88
$ws.Range("A88").Value = 1753876086
$ws.Range("B88").Value = "delete"
$ws.Range("C88").Value = "freq"
$ws.Range("D88").Value = "ser_pub_loc___variable_13---value_2 long with a lot of text"
$ws.Range("E88").Value = "ser_pub_loc___variable_13"

# Row 89
$ws.Range("A89").Value = 1753876146
$ws.Range("B89").Value = "add"
$ws.Range("C89").Value = "freq"
$ws.Range("D89").Value = "ser_pub_loc___variable_13---value_2 long with a lot of text, but like a lot of text very very very value_2 long with a lot of text, but like a lot of text very very very"

# Row 90
$ws.Range("A90").Value = 1753876146
$ws.Range("B90").Value = "delete"
$ws.Range("C90").Value = "freq"
$ws.Range("D90").Value = "ser_pub_loc___variable_13---value_2 long with a lot of text, but like a lot of text very very very"
$ws.Range("E90").Value = "ser_pub_loc___variable_13"

# Row 91
$ws.Range("A91").Value = 1753879819
$ws.Range("B91").Value = "add"
$ws.Range("C91").Value = "freq"
$ws.Range("D91").Value = "ser_pub_loc___variable_13---value_2"

# Row 92
$ws.Range("A92").Value = 1753879819
$ws.Range("B92").Value = "delete"
$ws.Range("C92").Value = "freq"
$ws.Range("D92").Value = "ser_pub_loc___variable_13---value_2 long with a lot of text, but like a lot of text very very very value_2 long with a lot of text, but like a lot of text very very very"
$ws.Range("E92").Value = "ser_pub_loc___variable_13"

